$d = $word.ActiveDocument

# 1. Replace the ID placeholder text (also removes the trailing-space run,
#    collapsing both runs into a single run that keeps the first run's
#    character formatting).
$r = $d.Range(0, 0)
$r.Find.Execute("**ID__AFFARS_mp_5315_3_topic_42__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP_5315_3_APPENDIX_B__ID**", 2) | Out-Null

# 2. Update the first paragraph's formatting: left indent 120 -> 225 twips
#    (6pt -> 11.25pt) and add a paragraph border with 5pt spacing on every
#    side.
$p = $d.Paragraphs.Item(1)
$p.Format.LeftIndent = 11.25
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5
